$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {

    # --- Unmerge the old "FREQUENCY REPORT" title block (was G4:H5) ---
    $ws.Range("G4:H5").UnMerge()

    # --- Row heights for the new 3-line header block in column G (rows 3 & 4) ---
    $ws.Rows.Item(3).RowHeight = 29
    $ws.Rows.Item(4).RowHeight = 29

    # --- G4: keep "FREQUENCY REPORT", now bold/size 22, vertical-center only ---
    $ws.Range("G4").Font.Size = 22
    $ws.Range("G4").HorizontalAlignment = 1

    # --- G3: new "client_name" placeholder, size 22, not bold ---
    $ws.Range("G3").Value = "client_name"
    $ws.Range("G3").Font.Size = 22

    # --- G5: new "date_time" placeholder, size 11, not bold, vertical-center ---
    $ws.Range("G5").Value = "date_time"
    $ws.Range("G5").Font.Bold = $false
    $ws.Range("G5").Font.Size = 11
    $ws.Range("G5").HorizontalAlignment = 1

    # --- Reposition/resize the logo picture (anchor moved up-left slightly) ---
    $shp = $ws.Shapes.Item(1)
    $shp.Left = 14
    $shp.Top = 10
    $shp.Width = 357
    $shp.Height = 120
}

# --- Active sheet / selection swap: "Current deliveries" becomes the active tab ---
$ws1 = $wb.Worksheets.Item("Current deliveries")
$ws2 = $wb.Worksheets.Item("Completed deliveries")

$ws2.Activate()
$ws2.Range("C25").Select()

$ws1.Activate()
$ws1.Range("G27").Select()

Write-Output "done"
